$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1386
$wsExhibit.Range("F3").Value = 2931
$wsExhibit.Range("F5").Value = 268

# Sheet "全部类型" - same events appear again, update matching rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1386
$wsAll.Range("F4").Value = 2931
$wsAll.Range("F7").Value = 268
